$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '301.11'
Set-TextValue $ws.Range('E2') '-4.61%'
Set-TextValue $ws.Range('D3') '35.14'
Set-TextValue $ws.Range('E3') '-1.26%'
Set-TextValue $ws.Range('E4') '-1.50%'
Set-TextValue $ws.Range('D5') '0.07937'
Set-TextValue $ws.Range('E5') '-2.06%'
Set-TextValue $ws.Range('D6') '1.911'
Set-TextValue $ws.Range('E6') '-10.51%'
Set-TextValue $ws.Range('D7') '7.800'
Set-TextValue $ws.Range('E7') '-2.52%'
Set-TextValue $ws.Range('B8') 'BTSEToken'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D8') '2.913'
Set-TextValue $ws.Range('E8') '2.90%'
Set-TextValue $ws.Range('B9') 'MXToken'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D9') '0.9213'
Set-TextValue $ws.Range('E9') '-0.58%'
Set-TextValue $ws.Range('B10') 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws.Range('C10') 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D10') '0.1314'
Set-TextValue $ws.Range('E10') '29.61%'
Set-TextValue $ws.Range('B11') 'WazirX'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D11') '0.1843'
Set-TextValue $ws.Range('E11') '-1.62%'
Set-TextValue $ws.Range('B12') 'MandalaExchangeToken'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D12') '0.09467'
Set-TextValue $ws.Range('E12') '2.61%'
Set-TextValue $ws.Range('B13') 'BitrueCoin'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D13') '0.03602'
Set-TextValue $ws.Range('E13') '-0.06%'
Set-TextValue $ws.Range('B14') 'BitMartToken'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D14') '0.09845'
Set-TextValue $ws.Range('E14') '-0.60%'
Set-TextValue $ws.Range('B15') 'BitForexToken'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D15') '0.001406'
Set-TextValue $ws.Range('E15') '-2.11%'
Set-TextValue $ws.Range('B16') 'TigerCash'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D16') '0.005807'
Set-TextValue $ws.Range('E16') '2.40%'
Set-TextValue $ws.Range('B17') 'LEO'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D17') '3.506'
Set-TextValue $ws.Range('E17') '0.82%'
Set-TextValue $ws.Range('B18') 'GateToken'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D18') '4.037'
Set-TextValue $ws.Range('E18') '-2.70%'
Set-TextValue $ws.Range('D19') '0.3425'
Set-TextValue $ws.Range('E19') '1.70%'
Set-TextValue $ws.Range('D20') '0.1312'
Set-TextValue $ws.Range('E20') '-1.40%'
Set-TextValue $ws.Range('D21') '5.042'
Set-TextValue $ws.Range('E21') '-2.18%'
Set-TextValue $ws.Range('D22') '0.2468'
Set-TextValue $ws.Range('E22') '11.14%'
Set-TextValue $ws.Range('D23') '0.04507'
Set-TextValue $ws.Range('E23') '-1.35%'
Set-TextValue $ws.Range('D24') '0.001215'
Set-TextValue $ws.Range('E24') '-2.55%'
Set-TextValue $ws.Range('D25') '0.004785'
Set-TextValue $ws.Range('E25') '1.66%'
Set-TextValue $ws.Range('E26') '0.00%'
Set-TextValue $ws.Range('D27') '0.0003008'
Set-TextValue $ws.Range('E27') '-33.22%'
Set-TextValue $ws.Range('D39') '0.01870'
Set-TextValue $ws.Range('E39') '-4.43%'
Set-TextValue $ws.Range('D40') '0.04702'
Set-TextValue $ws.Range('E40') '-3.43%'
Set-TextValue $ws.Range('D41') '0.007488'
Set-TextValue $ws.Range('E41') '-3.24%'
Set-TextValue $ws.Range('D42') '0.009768'
Set-TextValue $ws.Range('E42') '24.82%'
Set-TextValue $ws.Range('D43') '0.1322'
Set-TextValue $ws.Range('E43') '-4.97%'
Set-TextValue $ws.Range('D44') '0.002113'
Set-TextValue $ws.Range('E44') '-1.41%'
Set-TextValue $ws.Range('D45') '0.009608'
Set-TextValue $ws.Range('E45') '-17.45%'
Set-TextValue $ws.Range('D46') '0.00006230'
Set-TextValue $ws.Range('E46') '-4.68%'
Set-TextValue $ws.Range('D47') '0.00000000752'
Set-TextValue $ws.Range('E47') '0.17%'
Set-TextValue $ws.Range('E48') '75.15%'
Set-TextValue $ws.Range('E49') '-12.48%'
Set-TextValue $ws.Range('D50') '0.00002106'
Set-TextValue $ws.Range('E50') '0.17%'
Set-TextValue $ws.Range('D51') '0.0002006'
Set-TextValue $ws.Range('E51') '0.17%'
